$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write the new "Radio Calls for AWACS" table in columns T:Y, rows 1-14 ---
# Cells are written in this specific order so that new shared-string
# entries get created in the same order as in the target workbook.

$ws.Range("U3").Value = "Tower"
$ws.Range("U3").Font.Bold = $false

$ws.Range("T1").Value = "UHF"
$ws.Range("T1").Font.Bold = $true
$ws.Range("T1").HorizontalAlignment = -4108
$ws.Range("T1").VerticalAlignment = -4160

$ws.Range("V3").Value = "Ready for departure"

$ws.Range("Y12").Value = "CSA1"
$ws.Range("Y12").Font.Bold = $true

$ws.Range("U4").Value = "Departure"
$ws.Range("U4").Font.Bold = $false

$ws.Range("V4").Value = "Report Airborne"

$ws.Range("U12").Value = "AWACS"
$ws.Range("U12").Font.Bold = $false

$ws.Range("V12").Value = "Vector to nearest threat"

$ws.Range("Y13").Value = "CSA2"
$ws.Range("Y13").Font.Bold = $true

$ws.Range("V13").Value = "Request Picture"

$ws.Range("V14").Value = "Declare"

$ws.Range("Y14").Value = "CSA3"
$ws.Range("Y14").Font.Bold = $true

$ws.Range("U7").Value = "Approach"
$ws.Range("U7").Font.Bold = $false

$ws.Range("V7").Value = "Request Unrestricted Landing"

$ws.Range("Y7").Value = "TTT6"
$ws.Range("Y7").Font.Bold = $true

$ws.Range("Y4").Value = "TTTTT1"
$ws.Range("Y4").Font.Bold = $true

$ws.Range("V8").Value = "Inbound for Landing"

$ws.Range("V9").Value = "Report on final"

$ws.Range("Y9").Value = "TT4"
$ws.Range("Y9").Font.Bold = $true

$ws.Range("Y8").Value = "TT2"
$ws.Range("Y8").Font.Bold = $true

$ws.Range("Y3").Value = "TT1"
$ws.Range("Y3").Font.Bold = $true

$ws.Range("U2").Value = "Takeoff"
$ws.Range("U2").Font.Bold = $true

$ws.Range("U6").Value = "Landing"
$ws.Range("U6").Font.Bold = $true

$ws.Range("U11").Value = "Other"
$ws.Range("U11").Font.Bold = $true

# Remaining cells that reuse already-registered strings
$ws.Range("U8").Value = "Tower"
$ws.Range("U8").Font.Bold = $false

$ws.Range("U9").Value = "Tower"
$ws.Range("U9").Font.Bold = $false

$ws.Range("U13").Value = "AWACS"
$ws.Range("U13").Font.Bold = $false

$ws.Range("U14").Value = "AWACS"
$ws.Range("U14").Font.Bold = $false

# Numeric "note index" column T (bold, centered like existing column style)
$cellsT = @{ "T3" = 3; "T4" = 4; "T7" = 7; "T8" = 8; "T9" = 8; "T12" = 6; "T13" = 6; "T14" = 6 }
foreach ($addr in $cellsT.Keys) {
    $r = $ws.Range($addr)
    $r.Value = $cellsT[$addr]
    $r.Font.Bold = $true
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4160
}

# --- Column U width ---
$ws.Range("U1").ColumnWidth = 9.71

# --- Selection matches the post-edit UI state ---
$ws.Range("Q16").Select()

Write-Output "edit complete"
